$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# Update the time_taken column (F) values on the data sheet
$data.Range("F2").Value = "2021-10-05 14:34:52.220149"
$data.Range("F3").Value = "2021-10-05 14:34:52.220157"
$data.Range("F4").Value = "2021-10-05 14:34:52.220160"
$data.Range("F5").Value = "2021-10-05 14:34:52.220163"

# Add a new "metadata" worksheet placed right after "data"
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (bold/boxed style, matching the "data" sheet's header formatting)
$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$data.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)

# Data row
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Mirror movements"
$meta.Range("C2").Value = 3696

# data_version must be stored as text "1.0" (not the number 1)
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.0"
$meta.Range("D2").ClearFormats()

$meta.Range("E2").Value = "2021-09-02T02:08:23.638661Z"
$meta.Range("F2").Value = "2021-10-05 14:34:52.216450"
$meta.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3696/?format=json"

$data.Select()
